$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.279757499694824
$ws.Range("B1").Value = 3.246673345565796
$ws.Range("C1").Value = 5.883370399475098
$ws.Range("D1").Value = 1.761871218681335
$ws.Range("E1").Value = 1.033935904502869
